$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('E2').Value = '2026-02-19 21:18:51'
$ws.Range('I2').Value = '3.5 mm'
$ws.Range('E3').Value = '2026-02-19 21:18:54'
$ws.Range('I3').Value = '5.3 mm'
$ws.Range('E4').Value = '2026-02-19 21:18:57'
$ws.Range('H4').Value = "'56%"
$ws.Range('J4').Value = '1010.0 hPa'
$ws.Range('E5').Value = '2026-02-19 21:19:00'
$ws.Range('I5').Value = '7.7 mm'
$ws.Range('E6').Value = '2026-02-19 21:19:03'
$ws.Range('H6').Value = "'73%"
$ws.Range('J6').Value = '1010.1 hPa'
$ws.Range('E7').Value = '2026-02-19 21:19:05'
$ws.Range('J7').Value = '1011.2 hPa'
$ws.Range('E8').Value = '2026-02-19 21:19:08'
$ws.Range('J8').Value = '1010.9 hPa'
$ws.Range('E9').Value = '2026-02-19 21:19:10'
$ws.Range('K9').Value = '10.5 MJ/m2'
$ws.Range('O9').Value = '10.4 °C'
$ws.Range('E10').Value = '2026-02-19 21:19:13'
$ws.Range('H10').Value = "'71%"
$ws.Range('N10').Value = '4.0 °C 20:59 TU'
$ws.Range('O10').Value = '10.2 °C'
$ws.Range('E11').Value = '2026-02-19 21:19:16'
$ws.Range('H11').Value = "'63%"
$ws.Range('E12').Value = '2026-02-19 21:19:19'
$ws.Range('E13').Value = '2026-02-19 21:19:21'
$ws.Range('J13').Value = '1011.4 hPa'
$ws.Range('K13').Value = '13.6 MJ/m2'
$ws.Range('E14').Value = '2026-02-19 21:19:24'
$ws.Range('E15').Value = '2026-02-19 21:19:27'
$ws.Range('O15').Value = '9.8 °C'
$ws.Range('E16').Value = '2026-02-19 21:19:29'
$ws.Range('I16').Value = '10.2 mm'
$ws.Range('E17').Value = '2026-02-19 21:19:32'
$ws.Range('H17').Value = "'80%"
$ws.Range('E18').Value = '2026-02-19 21:19:35'
$ws.Range('H18').Value = "'60%"
$ws.Range('J18').Value = '1010.3 hPa'
$ws.Range('O18').Value = '11.6 °C'
$ws.Range('E19').Value = '2026-02-19 21:19:38'
$ws.Range('E20').Value = '2026-02-19 21:19:41'
$ws.Range('E21').Value = '2026-02-19 21:19:43'
$ws.Range('J21').Value = '1011.5 hPa'
$ws.Range('E22').Value = '2026-02-19 21:19:46'
$ws.Range('G22').Value = '144 cm'
$ws.Range('L22').Value = '101.5 km/h - 324º 20:57 TU'
$ws.Range('E23').Value = '2026-02-19 21:19:49'
$ws.Range('H23').Value = "'78%"
$ws.Range('I23').Value = '10.6 mm'
$ws.Range('E24').Value = '2026-02-19 21:19:51'
$ws.Range('J24').Value = '1015.0 hPa'
$ws.Range('E25').Value = '2026-02-19 21:19:54'
$ws.Range('I25').Value = '6.7 mm'
$ws.Range('E26').Value = '2026-02-19 21:19:57'
$ws.Range('H26').Value = "'56%"
$ws.Range('J26').Value = '1010.0 hPa'
$ws.Range('E27').Value = '2026-02-19 21:19:59'
$ws.Range('E28').Value = '2026-02-19 21:20:02'
$ws.Range('H28').Value = "'65%"
$ws.Range('J28').Value = '1010.0 hPa'
$ws.Range('O28').Value = '9.2 °C'
$ws.Range('E29').Value = '2026-02-19 21:20:05'
$ws.Range('N29').Value = '5.3 °C 20:59 TU'
$ws.Range('O29').Value = '10.4 °C'
$ws.Range('E30').Value = '2026-02-19 21:20:08'
$ws.Range('J30').Value = '1010.2 hPa'
$ws.Range('E31').Value = '2026-02-19 21:20:10'
$ws.Range('J31').Value = '1009.6 hPa'
$ws.Range('E32').Value = '2026-02-19 21:20:13'
$ws.Range('E33').Value = '2026-02-19 21:20:15'
$ws.Range('H33').Value = "'59%"
$ws.Range('J33').Value = '1010.9 hPa'
$ws.Range('O33').Value = '3.7 °C'
$ws.Range('E34').Value = '2026-02-19 21:20:18'
$ws.Range('E35').Value = '2026-02-19 21:20:21'
$ws.Range('H35').Value = "'68%"
$ws.Range('J35').Value = '1016.4 hPa'
$ws.Range('E36').Value = '2026-02-19 21:20:23'
$ws.Range('J36').Value = '1010.4 hPa'
$ws.Range('E37').Value = '2026-02-19 21:20:26'
$ws.Range('J37').Value = '1011.4 hPa'
$ws.Range('O37').Value = '5.9 °C'
$ws.Range('E38').Value = '2026-02-19 21:20:29'
$ws.Range('H38').Value = "'56%"
$ws.Range('K38').Value = '12.2 MJ/m2'
$ws.Range('L38').Value = '40.7 km/h - 271º 20:38 TU'
$ws.Range('E39').Value = '2026-02-19 21:20:31'
$ws.Range('E40').Value = '2026-02-19 21:20:34'
$ws.Range('J40').Value = '1012.7 hPa'
$ws.Range('O40').Value = '6.5 °C'
$ws.Range('E41').Value = '2026-02-19 21:20:37'
$ws.Range('J41').Value = '1013.1 hPa'
$ws.Range('E42').Value = '2026-02-19 21:20:40'
$ws.Range('H42').Value = "'76%"
$ws.Range('O42').Value = '11.1 °C'
$ws.Range('E43').Value = '2026-02-19 21:20:42'
$ws.Range('E44').Value = '2026-02-19 21:20:45'
$ws.Range('I44').Value = '9.2 mm'
$ws.Range('E45').Value = '2026-02-19 21:20:48'
$ws.Range('H45').Value = "'84%"
$ws.Range('J45').Value = '1015.7 hPa'
$ws.Range('E46').Value = '2026-02-19 21:20:51'
$ws.Range('J46').Value = '1015.9 hPa'
$ws.Range('O46').Value = '12.8 °C'
